# The workbook tracks daily price observations for "Haba" at Feria Lagunitas
# de Puerto Montt. This edit adds one new weekly observation, inserted as a
# new row right before the existing row 68 (pushing all following rows down
# by one, so the sheet grows from 171 to 172 data/header rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 68 - this shifts rows 68:171 down to 69:172
# and carries the row formatting (e.g. the date style on column D) down
# with them, matching how Excel's own "Insert Row" behaves.
$ws.Rows.Item(68).Insert()

# Populate the newly inserted row with the new observation's data.
$ws.Cells.Item(68, 1).Value  = 4
$ws.Cells.Item(68, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(68, 3).Value  = "Los Lagos"
$ws.Cells.Item(68, 4).Value  = 45219
$ws.Cells.Item(68, 5).Value  = 10
$ws.Cells.Item(68, 6).Value  = 100112026
$ws.Cells.Item(68, 7).Value  = "Haba"
$ws.Cells.Item(68, 8).Value  = "Sin especificar"
$ws.Cells.Item(68, 9).Value  = "Primera"
$ws.Cells.Item(68, 10).Value = 150
$ws.Cells.Item(68, 11).Value = 16000
$ws.Cells.Item(68, 12).Value = 16000
$ws.Cells.Item(68, 13).Value = 16000
$ws.Cells.Item(68, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(68, 15).Value = "Región Metropolitana"
$ws.Cells.Item(68, 16).Value = 640
$ws.Cells.Item(68, 17).Value = 25
$ws.Cells.Item(68, 18).Value = "Hortaliza"
